$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 125001496
$ws.Range("J43").Value = 1993
$ws.Range("L43").Value = 1993
$ws.Range("N43").Value = -2131
$ws.Range("H86").Value = 5364.697
$ws.Range("I86").Value = 5031.643
$ws.Range("K86").Value = 5031.643
$ws.Range("M86").Value = -3908.643
$ws.Range("H89").Value = 5364.697
$ws.Range("I89").Value = 5031.643
$ws.Range("K89").Value = 25158.215
$ws.Range("M89").Value = -19542.215
$ws.Range("H125").Value = 4764249
$ws.Range("I125").Value = 1203.8572
$ws.Range("J125").Value = 5293476.5
$ws.Range("K125").Value = 10834.7148
$ws.Range("L125").Value = 47641288.5
$ws.Range("M125").Value = -8374.7148
$ws.Range("N125").Value = -47646208.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3536263
$ws.Range("I2").Value = 3771814.2
$ws.Range("K2").Value = 3771814.2
$ws.Range("M2").Value = -3771701.2
$ws.Range("H32").Value = 3118.9507
$ws.Range("I32").Value = 2092.1594
$ws.Range("K32").Value = 2092.1594
$ws.Range("M32").Value = -1805.1594
$ws.Range("H45").Value = 7574602
$ws.Range("I45").Value = 11067355
$ws.Range("J45").Value = 6969.3335
$ws.Range("K45").Value = 11067355
$ws.Range("L45").Value = 6969.3335
$ws.Range("M45").Value = -11066978
$ws.Range("N45").Value = -7723.3335
$ws.Range("H61").Value = 8050.9375
$ws.Range("I61").Value = 9214.846
$ws.Range("J61").Value = 3007.3333
$ws.Range("K61").Value = 9214.846
$ws.Range("L61").Value = 3007.3333
$ws.Range("M61").Value = -9002.846
$ws.Range("N61").Value = -3431.3333
$ws.Range("H74").Value = 36840.965
$ws.Range("I74").Value = 7980.4824
$ws.Range("K74").Value = 7980.4824
$ws.Range("M74").Value = -7106.4824
$ws.Range("H77").Value = 36840.965
$ws.Range("I77").Value = 7980.4824
$ws.Range("K77").Value = 39902.412
$ws.Range("M77").Value = -35534.412
$ws.Range("H116").Value = 3536263
$ws.Range("I116").Value = 3771814.2
$ws.Range("K116").Value = 3771814.2
$ws.Range("M116").Value = -3769520.2
$ws.Range("H122").Value = 476296.2
$ws.Range("I122").Value = 2495.6453
$ws.Range("J122").Value = 1606128.2
$ws.Range("K122").Value = 7486.9359
$ws.Range("L122").Value = 4818384.6
$ws.Range("M122").Value = -5036.9359
$ws.Range("N122").Value = -4823284.6
$ws.Range("H132").Value = 4411.854
$ws.Range("I132").Value = 4418.375
$ws.Range("K132").Value = 13255.125
$ws.Range("M132").Value = -10725.125
$ws.Range("H136").Value = 8050.9375
$ws.Range("I136").Value = 9214.846
$ws.Range("J136").Value = 3007.3333
$ws.Range("K136").Value = 27644.538
$ws.Range("L136").Value = 9021.999899999999
$ws.Range("M136").Value = -25094.538
$ws.Range("N136").Value = -14121.9999
$ws.Range("H139").Value = 508063.9
$ws.Range("I139").Value = 500650
$ws.Range("J139").Value = 508887.66
$ws.Range("K139").Value = 500650
$ws.Range("L139").Value = 508887.66
$ws.Range("M139").Value = -495510
$ws.Range("N139").Value = -519167.66

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3536263
$ws.Range("I3").Value = 3771814.2
$ws.Range("K3").Value = 3771814.2
$ws.Range("M3").Value = -3771700.2
$ws.Range("H20").Value = 2819.353
$ws.Range("I20").Value = 2256.1538
$ws.Range("J20").Value = 4649.75
$ws.Range("K20").Value = 2256.1538
$ws.Range("L20").Value = 4649.75
$ws.Range("M20").Value = -2009.1538
$ws.Range("N20").Value = -5143.75
$ws.Range("H86").Value = 4770798
$ws.Range("I86").Value = 7154768.5
$ws.Range("J86").Value = 2857
$ws.Range("K86").Value = 7154768.5
$ws.Range("L86").Value = 2857
$ws.Range("M86").Value = -7153645.5
$ws.Range("N86").Value = -5103
$ws.Range("H89").Value = 4770798
$ws.Range("I89").Value = 7154768.5
$ws.Range("J89").Value = 2857
$ws.Range("K89").Value = 35773842.5
$ws.Range("L89").Value = 14285
$ws.Range("M89").Value = -35768226.5
$ws.Range("N89").Value = -25517
$ws.Range("H99").Value = 5718222
$ws.Range("I99").Value = 10992808
$ws.Range("J99").Value = 4086.9167
$ws.Range("K99").Value = 10992808
$ws.Range("L99").Value = 4086.9167
$ws.Range("M99").Value = -10991310
$ws.Range("N99").Value = -7082.9167
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 76841.5
$ws.Range("J9").Value = 76841.5
$ws.Range("L9").Value = 76841.5
$ws.Range("N9").Value = -77177.5
$ws.Range("H22").Value = 987.2
$ws.Range("I22").Value = 487
$ws.Range("K22").Value = 487
$ws.Range("M22").Value = -137
$ws.Range("H122").Value = 3573.35
$ws.Range("I122").Value = 2764.8572
$ws.Range("J122").Value = 4008.6924
$ws.Range("K122").Value = 8294.571599999999
$ws.Range("L122").Value = 12026.0772
$ws.Range("M122").Value = -5844.571599999999
$ws.Range("N122").Value = -16926.0772
$ws.Range("H132").Value = 60927.59
$ws.Range("I132").Value = 68651.664
$ws.Range("K132").Value = 205954.992
$ws.Range("M132").Value = -203424.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6029148
$ws.Range("I4").Value = 9402772
$ws.Range("J4").Value = 245791.14
$ws.Range("K4").Value = 28208316
$ws.Range("L4").Value = 737373.42
$ws.Range("M4").Value = -28208204
$ws.Range("N4").Value = -737597.42
$ws.Range("H12").Value = 83123.73
$ws.Range("I12").Value = 177791.4
$ws.Range("K12").Value = 533374.2
$ws.Range("M12").Value = -533201.2
$ws.Range("H87").Value = 20209.334
$ws.Range("I87").Value = 16404.8
$ws.Range("J87").Value = 21210.525
$ws.Range("K87").Value = 49214.39999999999
$ws.Range("L87").Value = 63631.575
$ws.Range("M87").Value = -47966.39999999999
$ws.Range("N87").Value = -66127.57500000001
$ws.Range("H90").Value = 20209.334
$ws.Range("I90").Value = 16404.8
$ws.Range("J90").Value = 21210.525
$ws.Range("K90").Value = 147643.2
$ws.Range("L90").Value = 190894.725
$ws.Range("M90").Value = -141403.2
$ws.Range("N90").Value = -203374.725
$ws.Range("H107").Value = 492.82608
$ws.Range("I107").Value = 372.1875
$ws.Range("J107").Value = 768.5714
$ws.Range("K107").Value = 1116.5625
$ws.Range("L107").Value = 2305.7142
$ws.Range("M107").Value = 803.4375
$ws.Range("N107").Value = -6145.7142
$ws.Range("H115").Value = 3395.2942
$ws.Range("I115").Value = 2520
$ws.Range("K115").Value = 7560
$ws.Range("M115").Value = -6385
$ws.Range("H129").Value = 1738.1482
$ws.Range("I129").Value = 1372.421
$ws.Range("J129").Value = 2606.75
$ws.Range("K129").Value = 4117.263
$ws.Range("L129").Value = 7820.25
$ws.Range("M129").Value = 882.7370000000001
$ws.Range("N129").Value = -17820.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20005658
$ws.Range("J70").Value = 6498.5
$ws.Range("L70").Value = 6498.5
$ws.Range("N70").Value = -7038.5
$ws.Range("H73").Value = 20005658
$ws.Range("J73").Value = 6498.5
$ws.Range("L73").Value = 6498.5
$ws.Range("N73").Value = -8370.5
$ws.Range("H132").Value = 5768.027
$ws.Range("I132").Value = 3194.3572
$ws.Range("J132").Value = 13775
$ws.Range("K132").Value = 9583.071599999999
$ws.Range("L132").Value = 41325
$ws.Range("M132").Value = -7053.071599999999
$ws.Range("N132").Value = -46385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3132461.2
$ws.Range("I40").Value = 5006374
$ws.Range("K40").Value = 5006374
$ws.Range("M40").Value = -5006238
$ws.Range("H46").Value = 4835517.5
$ws.Range("I46").Value = 14494020
$ws.Range("K46").Value = 14494020
$ws.Range("M46").Value = -14493832
$ws.Range("H93").Value = 30307848
$ws.Range("I93").Value = 33337984
$ws.Range("J93").Value = 6500
$ws.Range("K93").Value = 33337984
$ws.Range("L93").Value = 6500
$ws.Range("M93").Value = -33336736
$ws.Range("N93").Value = -8996
$ws.Range("H132").Value = 3381.7942
$ws.Range("I132").Value = 3005.261
$ws.Range("J132").Value = 4169.091
$ws.Range("K132").Value = 9015.782999999999
$ws.Range("L132").Value = 12507.273
$ws.Range("M132").Value = -6485.782999999999
$ws.Range("N132").Value = -17567.273
$ws.Range("H136").Value = 33823.03
$ws.Range("I136").Value = 44256.707
$ws.Range("J136").Value = 5999.8887
$ws.Range("K136").Value = 132770.121
$ws.Range("L136").Value = 17999.6661
$ws.Range("M136").Value = -130220.121
$ws.Range("N136").Value = -23099.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 152753
$ws.Range("J41").Value = 152753
$ws.Range("L41").Value = 152753
$ws.Range("N41").Value = -153533
$ws.Range("H62").Value = 19757
$ws.Range("J62").Value = 8293.643
$ws.Range("L62").Value = 8293.643
$ws.Range("N62").Value = -9541.643
$ws.Range("H65").Value = 19757
$ws.Range("J65").Value = 8293.643
$ws.Range("L65").Value = 41468.215
$ws.Range("N65").Value = -47708.215
$ws.Range("H126").Value = 4588
$ws.Range("I126").Value = 4076.0833
$ws.Range("K126").Value = 12228.2499
$ws.Range("M126").Value = -9758.249899999999
$ws.Range("H132").Value = 20636928
$ws.Range("I132").Value = 25004918
$ws.Range("K132").Value = 75014754
$ws.Range("M132").Value = -75012224
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
$ws.Range("H136").Value = 4614.396
$ws.Range("I136").Value = 5368.162
$ws.Range("K136").Value = 16104.486
$ws.Range("M136").Value = -13554.486
